# Apply updated natmi output values (Sirpa-Cd47) per Dr Hou's advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 79.465682
$ws.Range("H2").Value = 238.397046
$ws.Range("I2").Value = 0.2316778771755457
$ws.Range("J2").Value = 0.2316778771755458
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 57.478532
$ws.Range("N2").Value = 172.435596
$ws.Range("O2").Value = 0.2414676574042868
$ws.Range("P2").Value = 0.2414676574042868
$ws.Range("Q2").Value = 4567.570745738824
$ws.Range("R2").Value = 41108.13671164942
$ws.Range("S2").Value = 0.05594271427397711
$ws.Range("T2").Value = 0.05594271427397712

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 79.465682
$ws.Range("H3").Value = 238.397046
$ws.Range("I3").Value = 0.2316778771755457
$ws.Range("J3").Value = 0.2316778771755458
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 64.84043133333334
$ws.Range("N3").Value = 194.521294
$ws.Range("O3").Value = 0.272395040623924
$ws.Range("P3").Value = 0.2723950406239241
$ws.Range("Q3").Value = 5152.589097077504
$ws.Range("R3").Value = 46373.30187369753
$ws.Range("S3").Value = 0.06310790476489726
$ws.Range("T3").Value = 0.06310790476489728

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 79.465682
$ws.Range("H4").Value = 238.397046
$ws.Range("I4").Value = 0.2316778771755457
$ws.Range("J4").Value = 0.2316778771755458
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 85.31555666666667
$ws.Range("N4").Value = 255.94667
$ws.Range("O4").Value = 0.3584111648579104
$ws.Range("P4").Value = 0.3584111648579105
$ws.Range("Q4").Value = 6779.658895726313
$ws.Range("R4").Value = 61016.93006153682
$ws.Range("S4").Value = 0.08303593783029525
$ws.Range("T4").Value = 0.08303593783029527

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 79.465682
$ws.Range("H5").Value = 238.397046
$ws.Range("I5").Value = 0.2316778771755457
$ws.Range("J5").Value = 0.2316778771755458
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.40370266666666
$ws.Range("N5").Value = 91.211108
$ws.Range("O5").Value = 0.1277261371138787
$ws.Range("P5").Value = 0.1277261371138788
$ws.Range("Q5").Value = 2416.050967731885
$ws.Range("R5").Value = 21744.45870958697
$ws.Range("S5").Value = 0.02959132030637611
$ws.Range("T5").Value = 0.02959132030637612

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.22335033333333
$ws.Range("H6").Value = 63.670051
$ws.Range("I6").Value = 0.06187552447834749
$ws.Range("J6").Value = 0.06187552447834749
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 57.478532
$ws.Range("N6").Value = 172.435596
$ws.Range("O6").Value = 0.2414676574042868
$ws.Range("P6").Value = 0.2414676574042868
$ws.Range("Q6").Value = 1219.887021281711
$ws.Range("R6").Value = 10978.9831915354
$ws.Range("S6").Value = 0.01494093794644817
$ws.Range("T6").Value = 0.01494093794644817

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 21.22335033333333
$ws.Range("H7").Value = 63.670051
$ws.Range("I7").Value = 0.06187552447834749
$ws.Range("J7").Value = 0.06187552447834749
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 64.84043133333334
$ws.Range("N7").Value = 194.521294
$ws.Range("O7").Value = 0.272395040623924
$ws.Range("P7").Value = 0.2723950406239241
$ws.Range("Q7").Value = 1376.131189951777
$ws.Range("R7").Value = 12385.180709566
$ws.Range("S7").Value = 0.01685458600390607
$ws.Range("T7").Value = 0.01685458600390607

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.22335033333333
$ws.Range("H8").Value = 63.670051
$ws.Range("I8").Value = 0.06187552447834749
$ws.Range("J8").Value = 0.06187552447834749
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 85.31555666666667
$ws.Range("N8").Value = 255.94667
$ws.Range("O8").Value = 0.3584111648579104
$ws.Range("P8").Value = 0.3584111648579105
$ws.Range("Q8").Value = 1810.681948020019
$ws.Range("R8").Value = 16296.13753218017
$ws.Range("S8").Value = 0.02217687880447867
$ws.Range("T8").Value = 0.02217687880447868

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.22335033333333
$ws.Range("H9").Value = 63.670051
$ws.Range("I9").Value = 0.06187552447834749
$ws.Range("J9").Value = 0.06187552447834749
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.40370266666666
$ws.Range("N9").Value = 91.211108
$ws.Range("O9").Value = 0.1277261371138787
$ws.Range("P9").Value = 0.1277261371138788
$ws.Range("Q9").Value = 645.2684331251675
$ws.Range("R9").Value = 5807.415898126508
$ws.Range("S9").Value = 0.00790312172351457
$ws.Range("T9").Value = 0.007903121723514574

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 236.9265036666667
$ws.Range("H10").Value = 710.779511
$ws.Range("I10").Value = 0.6907463452728876
$ws.Range("J10").Value = 0.6907463452728876
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 57.478532
$ws.Range("N10").Value = 172.435596
$ws.Range("O10").Value = 0.2414676574042868
$ws.Range("P10").Value = 0.2414676574042868
$ws.Range("Q10").Value = 13618.18762265262
$ws.Range("R10").Value = 122563.6886038736
$ws.Range("S10").Value = 0.1667929018536168
$ws.Range("T10").Value = 0.1667929018536168

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 236.9265036666667
$ws.Range("H11").Value = 710.779511
$ws.Range("I11").Value = 0.6907463452728876
$ws.Range("J11").Value = 0.6907463452728876
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 64.84043133333334
$ws.Range("N11").Value = 194.521294
$ws.Range("O11").Value = 0.272395040623924
$ws.Range("P11").Value = 0.2723950406239241
$ws.Range("Q11").Value = 15362.41669204525
$ws.Range("R11").Value = 138261.7502284072
$ws.Range("S11").Value = 0.1881558787814352
$ws.Range("T11").Value = 0.1881558787814353

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 236.9265036666667
$ws.Range("H12").Value = 710.779511
$ws.Range("I12").Value = 0.6907463452728876
$ws.Range("J12").Value = 0.6907463452728876
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 85.31555666666667
$ws.Range("N12").Value = 255.94667
$ws.Range("O12").Value = 0.3584111648579104
$ws.Range("P12").Value = 0.3584111648579105
$ws.Range("Q12").Value = 20213.51654940871
$ws.Range("R12").Value = 181921.6489446783
$ws.Range("S12").Value = 0.2475712022306
$ws.Range("T12").Value = 0.2475712022306001

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 236.9265036666667
$ws.Range("H13").Value = 710.779511
$ws.Range("I13").Value = 0.6907463452728876
$ws.Range("J13").Value = 0.6907463452728876
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.40370266666666
$ws.Range("N13").Value = 91.211108
$ws.Range("O13").Value = 0.1277261371138787
$ws.Range("P13").Value = 0.1277261371138788
$ws.Range("Q13").Value = 7203.442971334242
$ws.Range("R13").Value = 64830.98674200818
$ws.Range("S13").Value = 0.08822636240723546
$ws.Range("T13").Value = 0.08822636240723548

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.385198333333332
$ws.Range("H14").Value = 16.155595
$ws.Range("I14").Value = 0.01570025307321912
$ws.Range("J14").Value = 0.01570025307321912
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 57.478532
$ws.Range("N14").Value = 172.435596
$ws.Range("O14").Value = 0.2414676574042868
$ws.Range("P14").Value = 0.2414676574042868
$ws.Range("Q14").Value = 309.5332947288466
$ws.Range("R14").Value = 2785.79965255962
$ws.Range("S14").Value = 0.003791103330244675
$ws.Range("T14").Value = 0.003791103330244676

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.385198333333332
$ws.Range("H15").Value = 16.155595
$ws.Range("I15").Value = 0.01570025307321912
$ws.Range("J15").Value = 0.01570025307321912
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 64.84043133333334
$ws.Range("N15").Value = 194.521294
$ws.Range("O15").Value = 0.272395040623924
$ws.Range("P15").Value = 0.2723950406239241
$ws.Range("Q15").Value = 349.1785827488811
$ws.Range("R15").Value = 3142.60724473993
$ws.Range("S15").Value = 0.00427667107368541
$ws.Range("T15").Value = 0.004276671073685411

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.385198333333332
$ws.Range("H16").Value = 16.155595
$ws.Range("I16").Value = 0.01570025307321912
$ws.Range("J16").Value = 0.01570025307321912
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 85.31555666666667
$ws.Range("N16").Value = 255.94667
$ws.Range("O16").Value = 0.3584111648579104
$ws.Range("P16").Value = 0.3584111648579105
$ws.Range("Q16").Value = 459.4411935687388
$ws.Range("R16").Value = 4134.970742118649
$ws.Range("S16").Value = 0.005627145992536453
$ws.Range("T16").Value = 0.005627145992536454

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 5.385198333333332
$ws.Range("H17").Value = 16.155595
$ws.Range("I17").Value = 0.01570025307321912
$ws.Range("J17").Value = 0.01570025307321912
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.40370266666666
$ws.Range("N17").Value = 91.211108
$ws.Range("O17").Value = 0.1277261371138787
$ws.Range("P17").Value = 0.1277261371138788
$ws.Range("Q17").Value = 163.7299689276955
$ws.Range("R17").Value = 1473.56972034926
$ws.Range("S17").Value = 0.002005332676752581
$ws.Range("T17").Value = 0.002005332676752582
